$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '22.238.90'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -1.04%  '
$ws.Range('E2').ClearFormats()
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.558.74'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -0.80%  '
$ws.Range('E3').ClearFormats()
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.28%  '
$ws.Range('E4').ClearFormats()
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  -0.25%  '
$ws.Range('E5').ClearFormats()
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '288.44'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('E6').ClearFormats()
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3814'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +2.32%  '
$ws.Range('E7').ClearFormats()
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3323'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -0.29%  '
$ws.Range('E8').ClearFormats()
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '44.80'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -7.40%  '
$ws.Range('E9').ClearFormats()
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.145'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +1.18%  '
$ws.Range('E10').ClearFormats()
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07410'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -1.04%  '
$ws.Range('E11').ClearFormats()
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.001'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -0.29%  '
$ws.Range('E12').ClearFormats()
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '20.24'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -3.22%  '
$ws.Range('E13').ClearFormats()
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.851'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -2.06%  '
$ws.Range('E14').ClearFormats()
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.754'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -2.36%  '
$ws.Range('E15').ClearFormats()
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '1.566.56'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -0.64%  '
$ws.Range('E16').ClearFormats()
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -3.77%  '
$ws.Range('E17').ClearFormats()
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.06661'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -1.47%  '
$ws.Range('E18').ClearFormats()
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '86.49'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -2.04%  '
$ws.Range('E19').ClearFormats()
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '6.412'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +0.27%  '
$ws.Range('E20').ClearFormats()
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.9995'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -0.41%  '
$ws.Range('E21').ClearFormats()
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '16.17'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -1.80%  '
$ws.Range('E22').ClearFormats()
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '11.75'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -2.82%  '
$ws.Range('E23').ClearFormats()
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '22.231.56'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -1.07%  '
$ws.Range('E24').ClearFormats()
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.278'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -4.63%  '
$ws.Range('E25').ClearFormats()
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.560'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -0.35%  '
$ws.Range('E26').ClearFormats()
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '151.30'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -0.91%  '
$ws.Range('E27').ClearFormats()
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '19.30'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -2.18%  '
$ws.Range('E28').ClearFormats()
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '4.939'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -1.33%  '
$ws.Range('E29').ClearFormats()
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '123.37'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -0.52%  '
$ws.Range('E30').ClearFormats()
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.736.19'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -0.92%  '
$ws.Range('E31').ClearFormats()
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.093'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +3.83%  '
$ws.Range('E32').ClearFormats()
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.919'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -3.88%  '
$ws.Range('E33').ClearFormats()
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.912'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -5.08%  '
$ws.Range('E34').ClearFormats()
$ws.Range('B35').NumberFormat = "@"
$ws.Range('B35').Value = 'Stellar'
$ws.Range('B35').ClearFormats()
$ws.Range('C35').NumberFormat = "@"
$ws.Range('C35').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('C35').ClearFormats()
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.08225'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -0.96%  '
$ws.Range('E35').ClearFormats()
$ws.Range('B36').NumberFormat = "@"
$ws.Range('B36').Value = 'FraxShare'
$ws.Range('B36').ClearFormats()
$ws.Range('C36').NumberFormat = "@"
$ws.Range('C36').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('C36').ClearFormats()
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '9.325'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -3.63%  '
$ws.Range('E36').ClearFormats()
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.06331'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -0.72%  '
$ws.Range('E37').ClearFormats()
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.02333'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -5.12%  '
$ws.Range('E38').ClearFormats()
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.323'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -1.00%  '
$ws.Range('E39').ClearFormats()
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.2165'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -4.69%  '
$ws.Range('E40').ClearFormats()
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.233'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -4.79%  '
$ws.Range('E41').ClearFormats()
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '11.01'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -2.57%  '
$ws.Range('E42').ClearFormats()
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.6072'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -3.67%  '
$ws.Range('E43').ClearFormats()
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -0.27%  '
$ws.Range('E44').ClearFormats()
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '13.78'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -0.89%  '
$ws.Range('E45').ClearFormats()
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -0.96%  '
$ws.Range('E46').ClearFormats()
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.5878'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -4.28%  '
$ws.Range('E47').ClearFormats()
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '122.32'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -2.20%  '
$ws.Range('E48').ClearFormats()
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -4.00%  '
$ws.Range('E49').ClearFormats()
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.179'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -2.84%  '
$ws.Range('E50').ClearFormats()
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.07056'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -2.90%  '
$ws.Range('E51').ClearFormats()
